$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.186.18'
$ws.Range("E2").Value = '  -1.73%  '

$ws.Range("D3").Value = '2.433.84'
$ws.Range("E3").Value = '  -2.15%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.11'
$ws.Range("E5").Value = '  -2.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.81'
$ws.Range("E6").Value = '  -2.41%  '

$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.498'
$ws.Range("E8").Value = '  -2.42%  '

$ws.Range("D9").Value = '2.431.26'
$ws.Range("E9").Value = '  -2.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.146'
$ws.Range("E10").Value = '  -7.63%  '

$ws.Range("E11").Value = '  -1.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.330'
$ws.Range("E12").Value = '  -6.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.74'
$ws.Range("E13").Value = '  -2.34%  '

$ws.Range("D14").Value = '2.887.46'
$ws.Range("E14").Value = '  -1.79%  '

$ws.Range("D15").Value = '68.125.80'
$ws.Range("E15").Value = '  -1.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000165'
$ws.Range("E16").Value = '  -4.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.07'
$ws.Range("E17").Value = '  -4.13%  '

$ws.Range("D18").Value = '2.414.92'
$ws.Range("E18").Value = '  -2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  -4.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '337.50'
$ws.Range("E20").Value = '  -2.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.91'
$ws.Range("E21").Value = '  -6.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.75'
$ws.Range("E22").Value = '  -3.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.83'
$ws.Range("E24").Value = '  -2.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.02'
$ws.Range("E25").Value = '  -4.61%  '

$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '2.572.06'
$ws.Range("E26").Value = '  -1.60%  '

$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.58'
$ws.Range("E27").Value = '  -7.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.95'
$ws.Range("E29").Value = '  -7.34%  '

$ws.Range("D30").Value = '0.0₃0796'
$ws.Range("E30").Value = '  -7.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.03'
$ws.Range("E31").Value = '  -6.61%  '

$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '430.00'
$ws.Range("E33").Value = '  -1.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.11'
$ws.Range("E34").Value = '  -6.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  -6.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.68'
$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.00'
$ws.Range("E37").Value = '  -0.25%  '

$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("E39").Value = '  -3.92%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.67'
$ws.Range("E40").Value = '  -2.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.298'
$ws.Range("E41").Value = '  -4.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.34'
$ws.Range("E42").Value = '  -4.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '37.33'
$ws.Range("E43").Value = '  -0.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.43'
$ws.Range("E44").Value = '  -8.81%  '

$ws.Range("E45").Value = '  +0.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.99'
$ws.Range("E46").Value = '  -7.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '130.26'
$ws.Range("E47").Value = '  -5.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.29'
$ws.Range("E48").Value = '  -3.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0709'
$ws.Range("E49").Value = '  -1.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.475'
$ws.Range("E50").Value = '  -5.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.552'
$ws.Range("E51").Value = '  -3.34%  '

